$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
# Row 6
$ws.Cells.Item(6, 8).Value = 123.375
$ws.Cells.Item(6, 9).Value = 136.71428
$ws.Cells.Item(6, 11).Value = 410.14284
$ws.Cells.Item(6, 13).Value = -298.14284

# Row 31
$ws.Cells.Item(31, 8).Value = 358.33334
$ws.Cells.Item(31, 9).Value = 358.33334
$ws.Cells.Item(31, 11).Value = 1075.00002
$ws.Cells.Item(31, 13).Value = -845.0000199999999

# Row 69
$ws.Cells.Item(69, 8).Value = 9945
$ws.Cells.Item(69, 10).Value = 9945
$ws.Cells.Item(69, 12).Value = 29835
$ws.Cells.Item(69, 14).Value = -31583

# Row 72
$ws.Cells.Item(72, 8).Value = 9945
$ws.Cells.Item(72, 10).Value = 9945
$ws.Cells.Item(72, 12).Value = 89505
$ws.Cells.Item(72, 14).Value = -98241

# Row 86
$ws.Cells.Item(86, 9).Value = 166668300
$ws.Cells.Item(86, 10).Value = 1969
$ws.Cells.Item(86, 11).Value = 166668300
$ws.Cells.Item(86, 12).Value = 1969
$ws.Cells.Item(86, 13).Value = -166667177
$ws.Cells.Item(86, 14).Value = -4215

# Row 88
$ws.Cells.Item(88, 8).Value = 1335.8667
$ws.Cells.Item(88, 9).Value = 1348.25
$ws.Cells.Item(88, 10).Value = 1331.3636
$ws.Cells.Item(88, 11).Value = 1348.25
$ws.Cells.Item(88, 12).Value = 1331.3636
$ws.Cells.Item(88, 13).Value = -942.25
$ws.Cells.Item(88, 14).Value = -2143.3636

# Row 89
$ws.Cells.Item(89, 9).Value = 166668300
$ws.Cells.Item(89, 10).Value = 1969
$ws.Cells.Item(89, 11).Value = 833341500
$ws.Cells.Item(89, 12).Value = 9845
$ws.Cells.Item(89, 13).Value = -833335884
$ws.Cells.Item(89, 14).Value = -21077

# Row 91
$ws.Cells.Item(91, 8).Value = 1335.8667
$ws.Cells.Item(91, 9).Value = 1348.25
$ws.Cells.Item(91, 10).Value = 1331.3636
$ws.Cells.Item(91, 11).Value = 1348.25
$ws.Cells.Item(91, 12).Value = 1331.3636
$ws.Cells.Item(91, 13).Value = 55.75
$ws.Cells.Item(91, 14).Value = -4139.3636

# Row 96
$ws.Cells.Item(96, 8).Value = 1796.1
$ws.Cells.Item(96, 9).Value = 1436.8889
$ws.Cells.Item(96, 10).Value = 5029
$ws.Cells.Item(96, 11).Value = 4310.6667
$ws.Cells.Item(96, 12).Value = 15087
$ws.Cells.Item(96, 13).Value = -2937.6667
$ws.Cells.Item(96, 14).Value = -17833

# Row 98
$ws.Cells.Item(98, 8).Value = 1816.1
$ws.Cells.Item(98, 9).Value = 1894.6666
$ws.Cells.Item(98, 11).Value = 1894.6666
$ws.Cells.Item(98, 13).Value = -396.6666

# Row 105
$ws.Cells.Item(105, 8).Value = 29832.334
$ws.Cells.Item(105, 10).Value = 29832.334
$ws.Cells.Item(105, 12).Value = 29832.334
$ws.Cells.Item(105, 14).Value = -36820.334

# Row 107
$ws.Cells.Item(107, 8).Value = 2348.8572
$ws.Cells.Item(107, 10).Value = 6728.6665
$ws.Cells.Item(107, 12).Value = 6728.6665
$ws.Cells.Item(107, 14).Value = -10568.6665

# Row 111
$ws.Cells.Item(111, 8).Value = 1776.3334
$ws.Cells.Item(111, 10).Value = 650
$ws.Cells.Item(111, 12).Value = 1950
$ws.Cells.Item(111, 14).Value = -8084

# Row 112
$ws.Cells.Item(112, 8).Value = 3262.1428
$ws.Cells.Item(112, 10).Value = 1869.4667
$ws.Cells.Item(112, 12).Value = 5608.4001
$ws.Cells.Item(112, 14).Value = -7824.4001

# Row 122
$ws.Cells.Item(122, 8).Value = 1816.1
$ws.Cells.Item(122, 9).Value = 1894.6666
$ws.Cells.Item(122, 11).Value = 5683.9998
$ws.Cells.Item(122, 13).Value = -3233.9998

# Row 137
$ws.Cells.Item(137, 8).Value = 2771.4285
$ws.Cells.Item(137, 9).Value = 1917.5
$ws.Cells.Item(137, 10).Value = 3910
$ws.Cells.Item(137, 11).Value = 5752.5
$ws.Cells.Item(137, 12).Value = 11730
$ws.Cells.Item(137, 13).Value = -3202.5
$ws.Cells.Item(137, 14).Value = -16830

# Row 138
$ws.Cells.Item(138, 8).Value = 3413.7317
$ws.Cells.Item(138, 9).Value = 951.1539
$ws.Cells.Item(138, 10).Value = 4557.0713
$ws.Cells.Item(138, 11).Value = 2853.4617
$ws.Cells.Item(138, 12).Value = 13671.2139
$ws.Cells.Item(138, 13).Value = 2286.5383
$ws.Cells.Item(138, 14).Value = -23951.2139

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
# Row 18
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(18, 12).ClearContents()

# Row 61
$ws.Cells.Item(61, 8).Value = 3582.3333
$ws.Cells.Item(61, 9).Value = 1012
$ws.Cells.Item(61, 10).Value = 3903.625
$ws.Cells.Item(61, 11).Value = 1012
$ws.Cells.Item(61, 12).Value = 3903.625
$ws.Cells.Item(61, 13).Value = -800
$ws.Cells.Item(61, 14).Value = -4327.625

# Row 74
$ws.Cells.Item(74, 8).Value = 2040.5
$ws.Cells.Item(74, 9).Value = 2206
$ws.Cells.Item(74, 11).Value = 2206
$ws.Cells.Item(74, 13).Value = -1332

# Row 77
$ws.Cells.Item(77, 8).Value = 2040.5
$ws.Cells.Item(77, 9).Value = 2206
$ws.Cells.Item(77, 11).Value = 11030
$ws.Cells.Item(77, 13).Value = -6662

# Row 110
$ws.Cells.Item(110, 8).Value = 4341.1816
$ws.Cells.Item(110, 9).Value = 2109.5
$ws.Cells.Item(110, 11).Value = 2109.5
$ws.Cells.Item(110, 13).Value = -64.5

# Row 122
$ws.Cells.Item(122, 8).Value = 4174.8184
$ws.Cells.Item(122, 9).Value = 1397.5
$ws.Cells.Item(122, 10).Value = 5761.857
$ws.Cells.Item(122, 11).Value = 4192.5
$ws.Cells.Item(122, 12).Value = 17285.571
$ws.Cells.Item(122, 13).Value = -1742.5
$ws.Cells.Item(122, 14).Value = -22185.571

# Row 136
$ws.Cells.Item(136, 8).Value = 3582.3333
$ws.Cells.Item(136, 9).Value = 1012
$ws.Cells.Item(136, 10).Value = 3903.625
$ws.Cells.Item(136, 11).Value = 3036
$ws.Cells.Item(136, 12).Value = 11710.875
$ws.Cells.Item(136, 13).Value = -486
$ws.Cells.Item(136, 14).Value = -16810.875

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
# Row 94
$ws.Cells.Item(94, 8).Value = 4815.381
$ws.Cells.Item(94, 9).Value = 804.2
$ws.Cells.Item(94, 10).Value = 8461.909
$ws.Cells.Item(94, 11).Value = 804.2
$ws.Cells.Item(94, 12).Value = 8461.909
$ws.Cells.Item(94, 13).Value = -353.2
$ws.Cells.Item(94, 14).Value = -9363.909

# Row 99
$ws.Cells.Item(99, 8).Value = 7937.048
$ws.Cells.Item(99, 9).Value = 7790.8647
$ws.Cells.Item(99, 11).Value = 7790.8647
$ws.Cells.Item(99, 13).Value = -6292.8647

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
# Row 16
$ws.Cells.Item(16, 8).Value = 9906.556
$ws.Cells.Item(16, 9).Value = 9322.714
$ws.Cells.Item(16, 11).Value = 9322.714
$ws.Cells.Item(16, 13).Value = -9035.714

# Row 31
$ws.Cells.Item(31, 8).Value = 3412.5
$ws.Cells.Item(31, 9).Value = 1400
$ws.Cells.Item(31, 10).Value = 3700
$ws.Cells.Item(31, 11).Value = 1400
$ws.Cells.Item(31, 12).Value = 3700
$ws.Cells.Item(31, 13).Value = -1105
$ws.Cells.Item(31, 14).Value = -4290

# Row 34
$ws.Cells.Item(34, 8).Value = 3412.5
$ws.Cells.Item(34, 9).Value = 1400
$ws.Cells.Item(34, 10).Value = 3700
$ws.Cells.Item(34, 11).Value = 1400
$ws.Cells.Item(34, 12).Value = 3700
$ws.Cells.Item(34, 13).Value = -1198
$ws.Cells.Item(34, 14).Value = -4104

# Row 58
$ws.Cells.Item(58, 8).Value = 15872.25
$ws.Cells.Item(58, 9).Value = 15349.5
$ws.Cells.Item(58, 11).Value = 15349.5
$ws.Cells.Item(58, 13).Value = -15146.5

# Row 99
$ws.Cells.Item(99, 8).Value = 6545.278
$ws.Cells.Item(99, 9).Value = 2785.7144
$ws.Cells.Item(99, 10).Value = 8937.727999999999
$ws.Cells.Item(99, 11).Value = 2785.7144
$ws.Cells.Item(99, 12).Value = 8937.727999999999
$ws.Cells.Item(99, 13).Value = -1287.7144
$ws.Cells.Item(99, 14).Value = -11933.728

# Row 105
$ws.Cells.Item(105, 8).Value = 4581.75
$ws.Cells.Item(105, 9).Value = 4216.4546
$ws.Cells.Item(105, 11).Value = 4216.4546
$ws.Cells.Item(105, 13).Value = -2469.4546

# Row 113
$ws.Cells.Item(113, 8).Value = 9906.556
$ws.Cells.Item(113, 9).Value = 9322.714
$ws.Cells.Item(113, 11).Value = 9322.714
$ws.Cells.Item(113, 13).Value = -7152.714

# Row 126
$ws.Cells.Item(126, 8).Value = 6545.278
$ws.Cells.Item(126, 9).Value = 2785.7144
$ws.Cells.Item(126, 10).Value = 8937.727999999999
$ws.Cells.Item(126, 11).Value = 8357.143199999999
$ws.Cells.Item(126, 12).Value = 26813.184
$ws.Cells.Item(126, 13).Value = -5887.143199999999
$ws.Cells.Item(126, 14).Value = -31753.184

# Row 134
$ws.Cells.Item(134, 8).Value = 7347.25
$ws.Cells.Item(134, 9).Value = 2294
$ws.Cells.Item(134, 10).Value = 10379.2
$ws.Cells.Item(134, 11).Value = 6882
$ws.Cells.Item(134, 12).Value = 31137.6
$ws.Cells.Item(134, 13).Value = -4347
$ws.Cells.Item(134, 14).Value = -36207.60000000001

# Row 136
$ws.Cells.Item(136, 8).Value = 15872.25
$ws.Cells.Item(136, 9).Value = 15349.5
$ws.Cells.Item(136, 11).Value = 46048.5
$ws.Cells.Item(136, 13).Value = -43498.5

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
# Row 2
$ws.Cells.Item(2, 8).Value = 115
$ws.Cells.Item(2, 9).Value = 55
$ws.Cells.Item(2, 10).Value = 175
$ws.Cells.Item(2, 11).Value = 330
$ws.Cells.Item(2, 12).Value = 1050
$ws.Cells.Item(2, 13).Value = -217
$ws.Cells.Item(2, 14).Value = -1276

# Row 52
$ws.Cells.Item(52, 8).Value = 1999.5
$ws.Cells.Item(52, 10).Value = 1999.5
$ws.Cells.Item(52, 12).Value = 5998.5
$ws.Cells.Item(52, 14).Value = -6530.5

# Row 81
$ws.Cells.Item(81, 8).Value = 2434.5
$ws.Cells.Item(81, 9).Value = 1940.3334
$ws.Cells.Item(81, 10).Value = 2928.6667
$ws.Cells.Item(81, 11).Value = 5821.0002
$ws.Cells.Item(81, 12).Value = 8786.000100000001
$ws.Cells.Item(81, 13).Value = -4698.0002
$ws.Cells.Item(81, 14).Value = -11032.0001

# Row 84
$ws.Cells.Item(84, 8).Value = 2434.5
$ws.Cells.Item(84, 9).Value = 1940.3334
$ws.Cells.Item(84, 10).Value = 2928.6667
$ws.Cells.Item(84, 11).Value = 17463.0006
$ws.Cells.Item(84, 12).Value = 26358.0003
$ws.Cells.Item(84, 13).Value = -11847.0006
$ws.Cells.Item(84, 14).Value = -37590.0003

# Row 129
$ws.Cells.Item(129, 8).Value = 41671172
$ws.Cells.Item(129, 10).Value = 83341820
$ws.Cells.Item(129, 12).Value = 250025460
$ws.Cells.Item(129, 14).Value = -250035460

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
# Row 99
$ws.Cells.Item(99, 8).Value = 3194.5
$ws.Cells.Item(99, 9).Value = 3194.5
$ws.Cells.Item(99, 11).Value = 3194.5
$ws.Cells.Item(99, 13).Value = -948.5

# Row 102
$ws.Cells.Item(102, 8).Value = 7382.8066
$ws.Cells.Item(102, 9).Value = 6596.8667
$ws.Cells.Item(102, 11).Value = 6596.8667
$ws.Cells.Item(102, 13).Value = -4974.8667

# Row 104
$ws.Cells.Item(104, 8).Value = 83917.75
$ws.Cells.Item(104, 10).Value = 81890.336
$ws.Cells.Item(104, 12).Value = 81890.336
$ws.Cells.Item(104, 14).Value = -88878.336

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
# Row 136
$ws.Cells.Item(136, 8).Value = 4484.091
$ws.Cells.Item(136, 9).Value = 2899
$ws.Cells.Item(136, 11).Value = 8697
$ws.Cells.Item(136, 13).Value = -6147

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
# Row 107
$ws.Cells.Item(107, 8).Value = 648.3333
$ws.Cells.Item(107, 9).Value = 478
$ws.Cells.Item(107, 11).Value = 1434
$ws.Cells.Item(107, 13).Value = 486

# Row 136
$ws.Cells.Item(136, 8).Value = 51958.523
$ws.Cells.Item(136, 9).Value = 93700
$ws.Cells.Item(136, 10).Value = 6042.9
$ws.Cells.Item(136, 11).Value = 281100
$ws.Cells.Item(136, 12).Value = 18128.7
$ws.Cells.Item(136, 13).Value = -278550
$ws.Cells.Item(136, 14).Value = -23228.7
